$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180) for every data row, from row 2 through row 351.
$ws.Range("C2:C351").Value = 45180
